$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.738.13'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.863.96'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.87%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.023'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.74%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("E6").Value = '  -0.75%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4373'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  +1.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07437'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8841'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.60'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.863.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.61%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.754'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.26%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.497'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.34%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07142'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '86.65'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.91%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009080'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.70%  '

$ws.Range("E19").Value = '  -0.75%  '

$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.740.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.293'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.89%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.49%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.093.20'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.044'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +6.50%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.998'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.11%  '

$ws.Range("E29").Value = '  +2.04%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.95'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.63%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09051'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.220'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.60%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7677'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.14%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.038'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.566'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.022'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.67%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.139'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.78%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01979'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05295'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.883'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.90%  '

$ws.Range("E41").Value = '  +1.11%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.945'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1680'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.685'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.77%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '109.96'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.714'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.56%  '

$ws.Range("E48").Value = '  -0.78%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06507'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4721'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.91%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.875'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.13%  '

